# "Generate Report for Handback" -- fills in the Latest Target File / Latest
# Handback File / Latest Handback DateTime columns on the per-locale sheets
# now that the de-de handback has completed, and flips the Overview/Status
# columns from "In Translation" to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$mdName1 = "4847067c-b6b7-40dc-bed2-c0ae93eca820.md"
$mdName2 = "aa70038e-b89c-49ec-94cf-e841776416e8.md"
$mdUrl1  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5bd30d8902b147afcc4cad9a55e4927d236303fc/e2e/$mdName1"
$mdUrl2  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5bd30d8902b147afcc4cad9a55e4927d236303fc/e2e/$mdName2"

$statusText = "Handed back: in sync with en-US"

# ---- Overview sheet: Status columns (zh-cn / de-de) for both rows ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $statusText
$zh.Range("C3").Value = $statusText

$zh.Hyperlinks.Add($zh.Range("I2"), $mdUrl1, "", "", $mdName1)
$zh.Range("J2").Value = "4847067c-b6b7-40dc-bed2-c0ae93eca820.c23aae31faef25faac0ba289957b0163e1d24c61.zh-cn.xlf"
$zh.Range("K2").Value = "2016-09-02 02:30:01"

$zh.Hyperlinks.Add($zh.Range("I3"), $mdUrl2, "", "", $mdName2)
$zh.Range("J3").Value = "aa70038e-b89c-49ec-94cf-e841776416e8.9e9a6d06e1e3975bb2660725475317bf8f59e0e2.zh-cn.xlf"
$zh.Range("K3").Value = "2016-09-02 02:30:01"

$zh.Columns.Item(3).ColumnWidth = 29.1
$zh.Columns.Item(9).ColumnWidth = 39.17
$zh.Columns.Item(10).ColumnWidth = 39.17

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $statusText
$de.Range("C3").Value = $statusText

$de.Hyperlinks.Add($de.Range("I2"), $mdUrl1, "", "", $mdName1)
$de.Range("J2").Value = "4847067c-b6b7-40dc-bed2-c0ae93eca820.c23aae31faef25faac0ba289957b0163e1d24c61.de-de.xlf"
$de.Range("K2").Value = "2016-09-02 02:30:24"

$de.Hyperlinks.Add($de.Range("I3"), $mdUrl2, "", "", $mdName2)
$de.Range("J3").Value = "aa70038e-b89c-49ec-94cf-e841776416e8.9e9a6d06e1e3975bb2660725475317bf8f59e0e2.de-de.xlf"
$de.Range("K3").Value = "2016-09-02 02:30:24"

$de.Columns.Item(3).ColumnWidth = 29.1
$de.Columns.Item(9).ColumnWidth = 39.17
$de.Columns.Item(10).ColumnWidth = 39.17

# ---- Overview sheet column widths (zh-cn / de-de) ----
$overview.Columns.Item(5).ColumnWidth = 29.1
$overview.Columns.Item(6).ColumnWidth = 29.1
